$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")
$data = @(
    @(1.02, 1.040570452553882, 1.040959909991495, 1.056001669731556, 1.062710247375454, 1.038062391104093, 1.045656343949804, 1.043740779275838, 1.058740578936548, 1.065430857243638, 1.047141295896671),
    @(1.02, 1.042117031338104, 1.042099111701669, 1.057370786840091, 1.064174928091394, 1.038466944016177, 1.04684561660383, 1.044689647409496, 1.059921833833354, 1.066708774600758, 1.048332257454279),
    @(1.02, 1.043115717716358, 1.042834287243415, 1.058255879167218, 1.065121683539132, 1.038726210118243, 1.047612750836656, 1.045301090708961, 1.06068483399401, 1.067534166658393, 1.049100481105706),
    @(1.02, 1.043535083063784, 1.043142890572905, 1.058627782356977, 1.065519469189068, 1.038834608138218, 1.047934685234723, 1.04555753931478, 1.061005281612355, 1.067880807249906, 1.049422872687545),
    @(1.02, 1.043605468328795, 1.043194679359421, 1.058690215598955, 1.065586245842094, 1.038852773684712, 1.047988706278321, 1.045600563009539, 1.061059067716068, 1.067938989140393, 1.049476970447221),
    @(1.02, 1.043121323180831, 1.042838412634184, 1.058260849293259, 1.065126999667026, 1.038727660882213, 1.047617054766536, 1.045304519747217, 1.060689117071843, 1.067538799871059, 1.04910479114766),
    @(1.02, 1.041093556849484, 1.041345318707092, 1.056464542384117, 1.063205452414565, 1.038199632286383, 1.046058766262127, 1.044061982445336, 1.059140072694251, 1.065863051008214, 1.047544289694837),
    @(1.02, 1.037504231644502, 1.038698994084504, 1.05329267628124, 1.059811541579408, 1.037249860200995, 1.043294136821569, 1.041852792461077, 1.056399856165995, 1.06289835068059, 1.044775734163366),
    @(1.02, 1.035099913121047, 1.036924123377699, 1.051173306978277, 1.057543187739036, 1.036603520484476, 1.041438025448738, 1.040366423414233, 1.05456555489886, 1.060913540547063, 1.042916986899486),
    @(1.02, 1.034055974077649, 1.036152975365341, 1.05025437219498, 1.05655950476916, 1.036320488872043, 1.040631121765205, 1.039719511908168, 1.053769424518002, 1.060052027577553, 1.042108937320169),
    @(1.02, 1.033667768644641, 1.035866136239592, 1.049912845874483, 1.056193891645151, 1.036214879845974, 1.040330913254938, 1.039478717031947, 1.053473418864202, 1.059731703456609, 1.041808302479384),
    @(1.02, 1.033751060121185, 1.035927682441935, 1.04998611328788, 1.056272327422642, 1.036237555047764, 1.0403953312451, 1.039530391227731, 1.053536926200437, 1.059800428667322, 1.04187281195048),
    @(1.02, 1.034023893936994, 1.036129273364063, 1.050226145504931, 1.05652928780484, 1.036311768977662, 1.040606316451987, 1.039699618031848, 1.053744962501899, 1.06002555605582, 1.042084096780562),
    @(1.02, 1.034191937212061, 1.036253426908501, 1.050374011468188, 1.056687578914449, 1.036357431123502, 1.040736246407102, 1.039803817376685, 1.053873102240034, 1.060164221853325, 1.042214211251116),
    @(1.02, 1.03516913483838, 1.036975246175384, 1.051234267005338, 1.057608439826271, 1.036622237446876, 1.041491508944742, 1.040409286615019, 1.05461835153782, 1.060970671900634, 1.04297054634819),
    @(1.02, 1.035781333233421, 1.037427318548499, 1.051773546848133, 1.05818567179626, 1.036787494450683, 1.041964403748933, 1.040788192168552, 1.055085322229107, 1.061475975138158, 1.043444112717244),
    @(1.02, 1.036138143672145, 1.037690752676124, 1.052087981124473, 1.058522220107185, 1.036883581185116, 1.042239927575491, 1.04100888304987, 1.055357518788026, 1.06177051008906, 1.043720027819239),
    @(1.02, 1.036259760587075, 1.037780534372184, 1.052195175212099, 1.058636950667628, 1.036916292667404, 1.04233382214291, 1.041084079148503, 1.055450300539107, 1.061870905227301, 1.043814055727713),
    @(1.02, 1.035715678617013, 1.037378841582969, 1.05171569954046, 1.058123754948889, 1.036769795497866, 1.041913698501568, 1.040747572162631, 1.055035239315646, 1.061421781613697, 1.043393335462612),
    @(1.02, 1.033943563415774, 1.036069920973723, 1.050155467389861, 1.056453625779995, 1.036289928060933, 1.040544200122479, 1.039649798899671, 1.053683708986462, 1.059959270543223, 1.042021892238741),
    @(1.02, 1.032826812920176, 1.03524463031776, 1.049173366886955, 1.055402217375336, 1.035985446295011, 1.039680311234883, 1.038956670871046, 1.052832281979271, 1.059037877354961, 1.041156776529835),
    @(1.02, 1.033419069130069, 1.035682354871282, 1.049694105736583, 1.05595971796983, 1.036147121551279, 1.040138546277932, 1.039324389681265, 1.053283799858561, 1.059526503712547, 1.041615662319205),
    @(1.02, 1.035745345960486, 1.037400747016002, 1.051741838615284, 1.058151732943978, 1.036777793834945, 1.041936610971221, 1.040765927578561, 1.055057870183426, 1.061446269952486, 1.0434162804706),
    @(1.02, 1.038434130769742, 1.039384983281265, 1.05411349537352, 1.060689929607728, 1.037497705730219, 1.044011123442376, 1.042426288440095, 1.057109561341642, 1.063666237932738, 1.045493738987414)
)

$rowStart = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $rowStart + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $colLetter = $cols[$j]
        $ws.Range("$colLetter$rowNum").Value = $rowVals[$j]
    }
}
